# Apply the price/volume updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.504.23'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '3.108.98'
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''598.07'
$ws.Range("E5").Value = '  -1.78%  '

$ws.Range("D6").Value = '''142.61'
$ws.Range("E6").Value = '  -2.33%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.107.01'
$ws.Range("E8").Value = '  -0.94%  '

$ws.Range("D9").Value = '''0.517'
$ws.Range("E9").Value = '  -1.03%  '

$ws.Range("E10").Value = '  -1.86%  '

$ws.Range("D11").Value = '''5.34'
$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("D12").Value = '''0.465'
$ws.Range("E12").Value = '  -1.19%  '

$ws.Range("D13").Value = '''0.0000251'
$ws.Range("E13").Value = '  -1.16%  '

$ws.Range("D14").Value = '''35.06'
$ws.Range("E14").Value = '  -0.89%  '

$ws.Range("D15").Value = '3.624.21'
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("E16").Value = '  +2.48%  '

$ws.Range("D17").Value = '63.625.10'
$ws.Range("E17").Value = '  -0.24%  '

$ws.Range("D18").Value = '3.105.61'
$ws.Range("E18").Value = '  -0.80%  '

$ws.Range("D19").Value = '''6.76'
$ws.Range("E19").Value = '  -1.53%  '

$ws.Range("D20").Value = '''480.13'
$ws.Range("E20").Value = '  +1.01%  '

$ws.Range("D21").Value = '''14.54'
$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").Value = '''0.703'
$ws.Range("E22").Value = '  -1.20%  '

$ws.Range("D23").Value = '''7.57'
$ws.Range("E23").Value = '  -4.95%  '

$ws.Range("D24").Value = '''87.03'
$ws.Range("E24").Value = '  +4.60%  '

$ws.Range("D25").Value = '''13.24'
$ws.Range("E25").Value = '  -3.15%  '

$ws.Range("E26").Value = '  -0.16%  '

$ws.Range("D27").Value = '''2.72'
$ws.Range("E27").Value = '  -2.77%  '

$ws.Range("D28").Value = '''8.21'

$ws.Range("D29").Value = '''7.04'
$ws.Range("E29").Value = '  -0.53%  '

$ws.Range("E30").Value = '  -2.42%  '

$ws.Range("D31").Value = '''27.04'
$ws.Range("E31").Value = '  +3.18%  '

$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("D33").Value = '''0.111'
$ws.Range("E33").Value = '  -8.76%  '

$ws.Range("E34").Value = '  -2.19%  '

$ws.Range("E35").Value = '  -2.16%  '

$ws.Range("D36").Value = '''5.98'
$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("E37").Value = '  -3.95%  '

$ws.Range("D38").Value = '''52.43'
$ws.Range("E38").Value = '  -0.32%  '

$ws.Range("E39").Value = '  -2.90%  '

$ws.Range("D40").Value = '''435.78'
$ws.Range("E40").Value = '  -5.02%  '

$ws.Range("E41").Value = '  -1.32%  '

$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").Value = '''8.25'
$ws.Range("E43").Value = '  -0.64%  '

$ws.Range("D44").Value = '2.852.82'
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("D45").Value = '''0.259'
$ws.Range("E45").Value = '  -3.17%  '

$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").Value = '''2.43'
$ws.Range("E46").Value = '  +0.94%  '

$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '''2.20'
$ws.Range("E47").Value = '  -4.14%  '

$ws.Range("E48").Value = '  +0.03%  '

$ws.Range("D49").Value = '''25.64'
$ws.Range("E49").Value = '  -2.73%  '

$ws.Range("D50").Value = '''0.113'
$ws.Range("E50").Value = '  -0.60%  '

$ws.Range("D51").Value = '''121.29'
$ws.Range("E51").Value = '  +1.91%  '
